$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Remove the four inline citation hyperlinks (the ones pointing to
# "previous post", "this paper" (x2), and "Nylund, Asparouhov, and Muthen
# (2007)"). Using Hyperlink.Delete() un-links the field but leaves the
# visible text as plain runs, which we then trim with Range deletes below.
# (Hyperlinks #1 "Young People Survey" and the References-section arxiv link
# are left in place; the latter is removed later along with its paragraph.)
# ---------------------------------------------------------------------------
$d.Hyperlinks.Item(5).Delete()
$d.Hyperlinks.Item(4).Delete()
$d.Hyperlinks.Item(3).Delete()
$d.Hyperlinks.Item(2).Delete()

# ---------------------------------------------------------------------------
# Step 2: "... Mahalanobis Distance (see my previous post on Mahalanobis for
# identifying outliers)." -> "... Mahalanobis Distance."
# ---------------------------------------------------------------------------
$r1 = $d.Content
[void]$r1.Find.Execute(" Distance (see my ")
$keepEnd = $r1.Start + 9   # length of " Distance"

$r2 = $d.Content
[void]$r2.Find.Execute("outliers).")
$delEnd = $r2.End

$d.Range($keepEnd, $delEnd).Text = "."

# ---------------------------------------------------------------------------
# Step 3: "... ellipsodial distribution (see Figure 2 from this paper for a
# visual). However, ..." -> "... ellipsodial distribution. However, ..."
# ---------------------------------------------------------------------------
$r1 = $d.Content
[void]$r1.Find.Execute(" distribution (see Figure 2 from ")
$keepEnd = $r1.Start + 13   # length of " distribution"

$r2 = $d.Content
[void]$r2.Find.Execute(" for a visual). However,")
$delEnd = $r2.Start + 14    # up to & including the closing paren

$d.Range($keepEnd, $delEnd).Delete()

# ---------------------------------------------------------------------------
# Step 4: "... (ICL) criterion. See this paper for more details. ICL isn't
# ..." -> "... (ICL) criterion. ICL isn't ..."
# ---------------------------------------------------------------------------
$r1 = $d.Content
[void]$r1.Find.Execute(" (ICL) criterion. See ")
$keepEnd = $r1.Start + 17   # length of " (ICL) criterion."

$r2 = $d.Content
[void]$r2.Find.Execute(" for more details. ICL isn")
$delEnd = $r2.Start + 18    # length of " for more details."

$d.Range($keepEnd, $delEnd).Delete()

# ---------------------------------------------------------------------------
# Step 5: "... Based on simulations by Nylund, Asparouhov, and Muthen (2007)
# BIC and BLRT ..." -> "... Based on simulations BIC and BLRT ..."
# ---------------------------------------------------------------------------
$r1 = $d.Content
[void]$r1.Find.Execute(" by Nylund, Asparouhov, and Muth")
$keepEnd = $r1.Start

$r2 = $d.Content
[void]$r2.Find.Execute(") BIC and BLRT")
$delEnd = $r2.Start + 1

$d.Range($keepEnd, $delEnd).Delete()

# ---------------------------------------------------------------------------
# Step 6: drop the whole "References & Resources" heading + its three
# reference entries + the trailing empty paragraph, so the document ends
# right after "...number of profiles." followed by the sectPr.
# ---------------------------------------------------------------------------
$total = $d.Paragraphs.Count
$headingIdx = 0
for ($i = 1; $i -le $total; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "References & Resources") {
        $headingIdx = $i
        break
    }
}

$pStart = $d.Paragraphs.Item($headingIdx)
$pEnd = $d.Paragraphs.Item($total)
$d.Range($pStart.Range.Start, $pEnd.Range.End).Delete()
